$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 16, shifting rows 16:86 down to 17:87
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with the new record
$ws.Cells.Item(16, 1).Value = 2
$ws.Cells.Item(16, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(16, 3).Value = "Coquimbo"
$ws.Cells.Item(16, 4).Value = 45071
$ws.Cells.Item(16, 5).Value = 4
$ws.Cells.Item(16, 6).Value = 100112026
$ws.Cells.Item(16, 7).Value = "Haba"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 500
$ws.Cells.Item(16, 11).Value = 10000
$ws.Cells.Item(16, 12).Value = 11000
$ws.Cells.Item(16, 13).Value = 10500
$ws.Cells.Item(16, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(16, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(16, 16).Value = 420
$ws.Cells.Item(16, 17).Value = 25
$ws.Cells.Item(16, 18).Value = "Hortaliza"
